$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "Förändrad" (C) column for existing rows 2-26: 45254 -> 45258 ---
for ($r = 2; $r -le 26; $r++) {
    $ws.Cells.Item($r, 3).Value = 45258
}

# --- 2. Row 26 becomes a non-last row, so it picks up an explicit row height ---
$ws.Rows.Item(26).RowHeight = 15

# --- 3. Append three new data rows (27, 28, 29) ---

# Row 27: A 59354-2023
$ws.Cells.Item(27, 1).Value = "A 59354-2023"
$ws.Cells.Item(27, 2).NumberFormat = $ws.Cells.Item(26, 2).NumberFormat
$ws.Cells.Item(27, 2).Value = 45253
$ws.Cells.Item(27, 3).NumberFormat = $ws.Cells.Item(26, 3).NumberFormat
$ws.Cells.Item(27, 3).Value = 45258
$ws.Cells.Item(27, 4).Value = "OKÄNT"
$ws.Cells.Item(27, 5).Value = "OKÄNT"
$ws.Cells.Item(27, 7).Value = 2.2
for ($col = 8; $col -le 17; $col++) {
    $ws.Cells.Item(27, $col).Value = 0
}
$ws.Cells.Item(27, 18).WrapText = $true
$ws.Rows.Item(27).RowHeight = 15

# Row 28: A 59637-2023
$ws.Cells.Item(28, 1).Value = "A 59637-2023"
$ws.Cells.Item(28, 2).NumberFormat = $ws.Cells.Item(26, 2).NumberFormat
$ws.Cells.Item(28, 2).Value = 45254
$ws.Cells.Item(28, 3).NumberFormat = $ws.Cells.Item(26, 3).NumberFormat
$ws.Cells.Item(28, 3).Value = 45258
$ws.Cells.Item(28, 4).Value = "OKÄNT"
$ws.Cells.Item(28, 5).Value = "OKÄNT"
$ws.Cells.Item(28, 6).Value = "Kommuner"
$ws.Cells.Item(28, 7).Value = 4.1
for ($col = 8; $col -le 17; $col++) {
    $ws.Cells.Item(28, $col).Value = 0
}
$ws.Cells.Item(28, 18).WrapText = $true
$ws.Rows.Item(28).RowHeight = 15

# Row 29: A 59686-2023 (last row - no explicit custom row height)
$ws.Cells.Item(29, 1).Value = "A 59686-2023"
$ws.Cells.Item(29, 2).NumberFormat = $ws.Cells.Item(26, 2).NumberFormat
$ws.Cells.Item(29, 2).Value = 45256
$ws.Cells.Item(29, 3).NumberFormat = $ws.Cells.Item(26, 3).NumberFormat
$ws.Cells.Item(29, 3).Value = 45258
$ws.Cells.Item(29, 4).Value = "OKÄNT"
$ws.Cells.Item(29, 5).Value = "OKÄNT"
$ws.Cells.Item(29, 7).Value = 1.6
for ($col = 8; $col -le 17; $col++) {
    $ws.Cells.Item(29, $col).Value = 0
}
$ws.Cells.Item(29, 18).WrapText = $true

Write-Output "done"
